$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# Ensure column A (Date) new cells remain text, not auto-converted to dates
$ws.Range("A78:A94").NumberFormat = "@"

$ws.Cells.Item(78, 1).Value = "2026-01-28"
$ws.Cells.Item(78, 2).Value = "17:24:09"
$ws.Cells.Item(78, 3).Value = "17:00"
$ws.Cells.Item(78, 4).Value = "Bedroom"
$ws.Cells.Item(78, 5).Value = "In Bed | HR=94 | BR=46"
$ws.Cells.Item(78, 6).Value = "Occupied"

$ws.Cells.Item(79, 1).Value = "2026-01-28"
$ws.Cells.Item(79, 2).Value = "17:24:10"
$ws.Cells.Item(79, 3).Value = "17:00"
$ws.Cells.Item(79, 4).Value = "Bedroom"
$ws.Cells.Item(79, 5).Value = "In Bed | HR=82 | BR=34"
$ws.Cells.Item(79, 6).Value = "Occupied"

$ws.Cells.Item(80, 1).Value = "2026-01-28"
$ws.Cells.Item(80, 2).Value = "17:24:10"
$ws.Cells.Item(80, 3).Value = "17:00"
$ws.Cells.Item(80, 4).Value = "Bedroom"
$ws.Cells.Item(80, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(80, 6).Value = "Occupied"

$ws.Cells.Item(81, 1).Value = "2026-01-28"
$ws.Cells.Item(81, 2).Value = "17:24:11"
$ws.Cells.Item(81, 3).Value = "17:00"
$ws.Cells.Item(81, 4).Value = "Bedroom"
$ws.Cells.Item(81, 5).Value = "In Bed | HR=52 | BR=4"
$ws.Cells.Item(81, 6).Value = "Occupied"

$ws.Cells.Item(82, 1).Value = "2026-01-28"
$ws.Cells.Item(82, 2).Value = "17:24:11"
$ws.Cells.Item(82, 3).Value = "17:00"
$ws.Cells.Item(82, 4).Value = "Bedroom"
$ws.Cells.Item(82, 5).Value = "In Bed | HR=90 | BR=42"
$ws.Cells.Item(82, 6).Value = "Occupied"

$ws.Cells.Item(83, 1).Value = "2026-01-28"
$ws.Cells.Item(83, 2).Value = "17:24:12"
$ws.Cells.Item(83, 3).Value = "17:00"
$ws.Cells.Item(83, 4).Value = "Bedroom"
$ws.Cells.Item(83, 5).Value = "In Bed | HR=56 | BR=8"
$ws.Cells.Item(83, 6).Value = "Occupied"

$ws.Cells.Item(84, 1).Value = "2026-01-28"
$ws.Cells.Item(84, 2).Value = "17:24:12"
$ws.Cells.Item(84, 3).Value = "17:00"
$ws.Cells.Item(84, 4).Value = "Bedroom"
$ws.Cells.Item(84, 5).Value = "In Bed | HR=58 | BR=10"
$ws.Cells.Item(84, 6).Value = "Occupied"

$ws.Cells.Item(85, 1).Value = "2026-01-28"
$ws.Cells.Item(85, 2).Value = "17:24:13"
$ws.Cells.Item(85, 3).Value = "17:00"
$ws.Cells.Item(85, 4).Value = "Bedroom"
$ws.Cells.Item(85, 5).Value = "In Bed | HR=107 | BR=59"
$ws.Cells.Item(85, 6).Value = "Occupied"

$ws.Cells.Item(86, 1).Value = "2026-01-28"
$ws.Cells.Item(86, 2).Value = "17:24:13"
$ws.Cells.Item(86, 3).Value = "17:00"
$ws.Cells.Item(86, 4).Value = "Bedroom"
$ws.Cells.Item(86, 5).Value = "In Bed | HR=57 | BR=9"
$ws.Cells.Item(86, 6).Value = "Occupied"

$ws.Cells.Item(87, 1).Value = "2026-01-28"
$ws.Cells.Item(87, 2).Value = "17:24:14"
$ws.Cells.Item(87, 3).Value = "17:00"
$ws.Cells.Item(87, 4).Value = "Bedroom"
$ws.Cells.Item(87, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(87, 6).Value = "Occupied"

$ws.Cells.Item(88, 1).Value = "2026-01-28"
$ws.Cells.Item(88, 2).Value = "17:24:15"
$ws.Cells.Item(88, 3).Value = "17:00"
$ws.Cells.Item(88, 4).Value = "Bedroom"
$ws.Cells.Item(88, 5).Value = "In Bed | HR=54 | BR=6"
$ws.Cells.Item(88, 6).Value = "Occupied"

$ws.Cells.Item(89, 1).Value = "2026-01-28"
$ws.Cells.Item(89, 2).Value = "17:24:16"
$ws.Cells.Item(89, 3).Value = "17:00"
$ws.Cells.Item(89, 4).Value = "Bedroom"
$ws.Cells.Item(89, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(89, 6).Value = "Occupied"

$ws.Cells.Item(90, 1).Value = "2026-01-28"
$ws.Cells.Item(90, 2).Value = "17:24:17"
$ws.Cells.Item(90, 3).Value = "17:00"
$ws.Cells.Item(90, 4).Value = "Bedroom"
$ws.Cells.Item(90, 5).Value = "In Bed | HR=55 | BR=7"
$ws.Cells.Item(90, 6).Value = "Occupied"

$ws.Cells.Item(91, 1).Value = "2026-01-28"
$ws.Cells.Item(91, 2).Value = "17:24:18"
$ws.Cells.Item(91, 3).Value = "17:00"
$ws.Cells.Item(91, 4).Value = "Bedroom"
$ws.Cells.Item(91, 5).Value = "In Bed | HR=75 | BR=27"
$ws.Cells.Item(91, 6).Value = "Occupied"

$ws.Cells.Item(92, 1).Value = "2026-01-28"
$ws.Cells.Item(92, 2).Value = "17:24:19"
$ws.Cells.Item(92, 3).Value = "17:00"
$ws.Cells.Item(92, 4).Value = "Bedroom"
$ws.Cells.Item(92, 5).Value = "In Bed | HR=91 | BR=43"
$ws.Cells.Item(92, 6).Value = "Occupied"

$ws.Cells.Item(93, 1).Value = "2026-01-28"
$ws.Cells.Item(93, 2).Value = "17:24:20"
$ws.Cells.Item(93, 3).Value = "17:00"
$ws.Cells.Item(93, 4).Value = "Bedroom"
$ws.Cells.Item(93, 5).Value = "In Bed | HR=50 | BR=2"
$ws.Cells.Item(93, 6).Value = "Occupied"

$ws.Cells.Item(94, 1).Value = "2026-01-28"
$ws.Cells.Item(94, 2).Value = "17:24:24"
$ws.Cells.Item(94, 3).Value = "17:00"
$ws.Cells.Item(94, 4).Value = "Bedroom"
$ws.Cells.Item(94, 5).Value = "In Bed | HR=49 | BR=1"
$ws.Cells.Item(94, 6).Value = "Occupied"

